# BP-813: Affiliate Mapping for True Independent Stations
# Rename header columns H1 and J1, and align their formatting (borders)
# with the rest of the header/data range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Affiliation Mismatch Note" -> "IsTrueIND"
$ws.Range("H1").Value = "IsTrueIND"

# Rename "SalesGroupName" -> "RepFirm"
$ws.Range("J1").Value = "RepFirm"

# I1:J1 (header) should carry the same bordered/bold style as the rest of
# row 1 (e.g. H1) instead of the borderless header style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# I2:J5 (data) should carry the same bordered style as the rest of the
# data rows (e.g. H2) instead of the borderless style.
$ws.Range("H2").Copy()
$ws.Range("I2:J5").PasteSpecial(-4122)

# Reset the active selection back to the top-left cell.
$ws.Range("A1").Select() | Out-Null
